# Sync attendance_reports: swap the first two comma-separated entries in the
# "Recorded By" column (G) of the "Session Analysis Results" sheet, for every
# data row, leaving any further entries (and the single exact value
# "admin@admin.com, System") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 157
$col = 7  # column G = "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "admin@admin.com, System") {
        $parts = $val -split ", "
        if ($parts.Length -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
